$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# New KEY/VALUE rows to append starting at row 23, as loaded from
# columns_number_parameters.xlsx / the country-specific parameters.xlsx
$rows = @(
    @("TAXDB_REGIMES", 6),
    @("MIN_START_YEAR", 2011),
    @("MAX_START_YEAR", 2020),
    @("MIN_START_YEAR_TRAINING", 2019),
    @("MAX_START_YEAR_TRAINING", 2019),
    @("MIN_CAPITAL_INCOME_PER_MONTH", 0),
    @("MAX_CAPITAL_INCOME_PER_MONTH", 4000),
    @("MIN_PERSONAL_PENSION_PER_MONTH", 0),
    @("MAX_PERSONAL_PENSION_PER_MONTH", 15000),
    @("MAX_CHILD_AGE_FOR_FORMAL_CARE", 14),
    @("MIN_AGE_MATERNITY", 18),
    @("MAX_AGE_MATERNITY", 44),
    @("BASE_PRICE_YEAR", 2015),
    @("PROB_NEWBORN_IS_MALE", 0.5)
)

$startRow = 23
$endRow = $startRow + $rows.Count - 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Apply the KEY (bold) / VALUE (regular) Helvetica Neue 10pt formatting used for
# these new parameter rows. Set Size before Name to keep the style table compact,
# then propagate the finished look with PasteSpecial so every row reuses the same
# two cell formats instead of minting a fresh one per row.
$keyCell = $ws.Range("A" + $startRow)
$keyCell.Font.Size = 10
$keyCell.Font.Name = "Helvetica Neue"
$keyCell.Font.Bold = $true

$valueCell = $ws.Range("B" + $startRow)
$valueCell.Font.Size = 10
$valueCell.Font.Name = "Helvetica Neue"

if ($endRow -gt $startRow) {
    $keyCell.Copy()
    $ws.Range("A" + ($startRow + 1) + ":A" + $endRow).PasteSpecial(-4122)
    $valueCell.Copy()
    $ws.Range("B" + ($startRow + 1) + ":B" + $endRow).PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}

$ws.Range("E24").Select()
